$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update existing parameter values on row 2
$ws.Range("C2").Value = 11
$ws.Range("F2").Value = 0.0005
$ws.Range("H2").Value = 10

# Add new transformer parameter columns
$ws.Range("K1").Value = "d_model"
$ws.Range("L1").Value = "num_layers"
$ws.Range("K2").Value = 16
$ws.Range("L2").Value = 1

# Update selection to match the saved state
$ws.Range("H7").Select()
